$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing columns C..L to D..M
$ws.Columns.Item(3).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 3).Value = "statut_name"

# New values for the inserted column (rows 2-4)
$ws.Cells.Item(2, 3).Value = "pas de résultat ni de publication"
$ws.Cells.Item(3, 3).Value = "résultat et / ou publication posté dans les 12 mois"
$ws.Cells.Item(4, 3).Value = "résultat et / ou publication posté dans les 12 mois"
